$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.187539
$ws.Range("H2").Value = 0.562617
$ws.Range("I2").Value = 0.002165158903160718
$ws.Range("J2").Value = 0.002165158903160718
$ws.Range("M2").Value = 4.443921666666667
$ws.Range("N2").Value = 13.331765
$ws.Range("O2").Value = 0.5356711567217328
$ws.Range("P2").Value = 0.5356711567217328
$ws.Range("Q2").Value = 0.8334086254450002
$ws.Range("R2").Value = 7.500677629005001
$ws.Range("S2").Value = 0.00115981317414246
$ws.Range("T2").Value = 0.00115981317414246

# Row 3
$ws.Range("G3").Value = 0.187539
$ws.Range("H3").Value = 0.562617
$ws.Range("I3").Value = 0.002165158903160718
$ws.Range("J3").Value = 0.002165158903160718
$ws.Range("O3").Value = 0.07358962454342892
$ws.Range("P3").Value = 0.07358962454342892
$ws.Range("Q3").Value = 0.114492309448
$ws.Range("R3").Value = 1.030430785032
$ws.Range("S3").Value = 0.0001593332307604596
$ws.Range("T3").Value = 0.0001593332307604596

# Row 4
$ws.Range("G4").Value = 0.187539
$ws.Range("H4").Value = 0.562617
$ws.Range("I4").Value = 0.002165158903160718
$ws.Range("J4").Value = 0.002165158903160718
$ws.Range("M4").Value = 0.071866
$ws.Range("N4").Value = 0.215598
$ws.Range("O4").Value = 0.008662741208451554
$ws.Range("P4").Value = 0.008662741208451554
$ws.Range("Q4").Value = 0.013477677774
$ws.Range("R4").Value = 0.121299099966
$ws.Range("S4").Value = 0.00001875621125325612
$ws.Range("T4").Value = 0.00001875621125325612

# Row 5
$ws.Range("G5").Value = 0.187539
$ws.Range("H5").Value = 0.562617
$ws.Range("I5").Value = 0.002165158903160718
$ws.Range("J5").Value = 0.002165158903160718
$ws.Range("M5").Value = 3.169702
$ws.Range("N5").Value = 9.509106000000001
$ws.Range("O5").Value = 0.3820764775263868
$ws.Range("P5").Value = 0.3820764775263868
$ws.Range("Q5").Value = 0.5944427433780001
$ws.Range("R5").Value = 5.349984690402001
$ws.Range("S5").Value = 0.0008272562870045422
$ws.Range("T5").Value = 0.0008272562870045424

# Row 6
$ws.Range("I6").Value = 0.739596852820732
$ws.Range("J6").Value = 0.7395968528207321
$ws.Range("M6").Value = 4.443921666666667
$ws.Range("N6").Value = 13.331765
$ws.Range("O6").Value = 0.5356711567217328
$ws.Range("P6").Value = 0.5356711567217328
$ws.Range("Q6").Value = 284.6841382371373
$ws.Range("R6").Value = 2562.157244134235
$ws.Range("S6").Value = 0.3961807016582347
$ws.Range("T6").Value = 0.3961807016582347

# Row 7
$ws.Range("I7").Value = 0.739596852820732
$ws.Range("J7").Value = 0.7395968528207321
$ws.Range("O7").Value = 0.07358962454342892
$ws.Range("P7").Value = 0.07358962454342892
$ws.Range("S7").Value = 0.05442665471257932
$ws.Range("T7").Value = 0.05442665471257933

# Row 8
$ws.Range("I8").Value = 0.739596852820732
$ws.Range("J8").Value = 0.7395968528207321
$ws.Range("M8").Value = 0.071866
$ws.Range("N8").Value = 0.215598
$ws.Range("O8").Value = 0.008662741208451554
$ws.Range("P8").Value = 0.008662741208451554
$ws.Range("Q8").Value = 4.603841339511334
$ws.Range("R8").Value = 41.434572055602
$ws.Range("S8").Value = 0.006406936134571234
$ws.Range("T8").Value = 0.006406936134571235

# Row 9
$ws.Range("I9").Value = 0.739596852820732
$ws.Range("J9").Value = 0.7395968528207321
$ws.Range("M9").Value = 3.169702
$ws.Range("N9").Value = 9.509106000000001
$ws.Range("O9").Value = 0.3820764775263868
$ws.Range("P9").Value = 0.3820764775263868
$ws.Range("Q9").Value = 203.0557579596994
$ws.Range("R9").Value = 1827.501821637294
$ws.Range("S9").Value = 0.2825825603153468
$ws.Range("T9").Value = 0.2825825603153468

# Row 10
$ws.Range("G10").Value = 19.62095333333333
$ws.Range("H10").Value = 58.86286
$ws.Range("I10").Value = 0.2265261188241786
$ws.Range("J10").Value = 0.2265261188241786
$ws.Range("M10").Value = 4.443921666666667
$ws.Range("N10").Value = 13.331765
$ws.Range("O10").Value = 0.5356711567217328
$ws.Range("P10").Value = 0.5356711567217328
$ws.Range("Q10").Value = 87.19397963865556
$ws.Range("R10").Value = 784.7458167479
$ws.Range("S10").Value = 0.1213435080982324
$ws.Range("T10").Value = 0.1213435080982324

# Row 11
$ws.Range("G11").Value = 19.62095333333333
$ws.Range("H11").Value = 58.86286
$ws.Range("I11").Value = 0.2265261188241786
$ws.Range("J11").Value = 0.2265261188241786
$ws.Range("O11").Value = 0.07358962454342892
$ws.Range("P11").Value = 0.07358962454342892
$ws.Range("Q11").Value = 11.97856584872889
$ws.Range("R11").Value = 107.80709263856
$ws.Range("S11").Value = 0.01666997203355147
$ws.Range("T11").Value = 0.01666997203355147

# Row 12
$ws.Range("G12").Value = 19.62095333333333
$ws.Range("H12").Value = 58.86286
$ws.Range("I12").Value = 0.2265261188241786
$ws.Range("J12").Value = 0.2265261188241786
$ws.Range("M12").Value = 0.071866
$ws.Range("N12").Value = 0.215598
$ws.Range("O12").Value = 0.008662741208451554
$ws.Range("P12").Value = 0.008662741208451554
$ws.Range("Q12").Value = 1.410079432253333
$ws.Range("R12").Value = 12.69071489028
$ws.Range("S12").Value = 0.001962337144328805
$ws.Range("T12").Value = 0.001962337144328805

# Row 13
$ws.Range("G13").Value = 19.62095333333333
$ws.Range("H13").Value = 58.86286
$ws.Range("I13").Value = 0.2265261188241786
$ws.Range("J13").Value = 0.2265261188241786
$ws.Range("M13").Value = 3.169702
$ws.Range("N13").Value = 9.509106000000001
$ws.Range("O13").Value = 0.3820764775263868
$ws.Range("P13").Value = 0.3820764775263868
$ws.Range("Q13").Value = 62.19257502257334
$ws.Range("R13").Value = 559.73317520316
$ws.Range("S13").Value = 0.08655030154806588
$ws.Range("T13").Value = 0.08655030154806589

# Row 14
$ws.Range("G14").Value = 2.746778666666666
$ws.Range("H14").Value = 8.240335999999999
$ws.Range("I14").Value = 0.03171186945192871
$ws.Range("J14").Value = 0.03171186945192871
$ws.Range("M14").Value = 4.443921666666667
$ws.Range("N14").Value = 13.331765
$ws.Range("O14").Value = 0.5356711567217328
$ws.Range("P14").Value = 0.5356711567217328
$ws.Range("Q14").Value = 12.20646923033778
$ws.Range("R14").Value = 109.85822307304
$ws.Range("S14").Value = 0.01698713379112323
$ws.Range("T14").Value = 0.01698713379112324

# Row 15
$ws.Range("G15").Value = 2.746778666666666
$ws.Range("H15").Value = 8.240335999999999
$ws.Range("I15").Value = 0.03171186945192871
$ws.Range("J15").Value = 0.03171186945192871
$ws.Range("O15").Value = 0.07358962454342892
$ws.Range("P15").Value = 0.07358962454342892
$ws.Range("Q15").Value = 1.676904713628444
$ws.Range("R15").Value = 15.092142422656
$ws.Range("S15").Value = 0.002333664566537667
$ws.Range("T15").Value = 0.002333664566537667

# Row 16
$ws.Range("G16").Value = 2.746778666666666
$ws.Range("H16").Value = 8.240335999999999
$ws.Range("I16").Value = 0.03171186945192871
$ws.Range("J16").Value = 0.03171186945192871
$ws.Range("M16").Value = 0.071866
$ws.Range("N16").Value = 0.215598
$ws.Range("O16").Value = 0.008662741208451554
$ws.Range("P16").Value = 0.008662741208451554
$ws.Range("Q16").Value = 0.1973999956586666
$ws.Range("R16").Value = 1.776599960928
$ws.Range("S16").Value = 0.0002747117182982588
$ws.Range("T16").Value = 0.0002747117182982589

# Row 17
$ws.Range("G17").Value = 2.746778666666666
$ws.Range("H17").Value = 8.240335999999999
$ws.Range("I17").Value = 0.03171186945192871
$ws.Range("J17").Value = 0.03171186945192871
$ws.Range("M17").Value = 3.169702
$ws.Range("N17").Value = 9.509106000000001
$ws.Range("O17").Value = 0.3820764775263868
$ws.Range("P17").Value = 0.3820764775263868
$ws.Range("Q17").Value = 8.706469833290667
$ws.Range("R17").Value = 78.358228499616
$ws.Range("S17").Value = 0.01211635937596955
$ws.Range("T17").Value = 0.01211635937596955
